$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (JudgeBotOpinions) updates
$ws.Range("C2").Value = "MSG: None`n`nMSG: The decision has been recorded as no decision regarding the movie for Friday.`n"
$ws.Range("C3").Value = "MSG: None`n`nMSG: The decision to acquire the rights for both movies has been recorded.`n"
$ws.Range("C4").Value = "MSG: None`n`nMSG: The decision has been recorded to acquire the rights for `"Barbie.`"`n"
$ws.Range("C5").Value = "MSG: None`n`nMSG: The committee could not reach a decision regarding the movie to be shown on Friday.`n"
$ws.Range("C6").Value = "MSG: None`n`nMSG: The decision has been recorded, and no movie will be acquired for Friday based on the conversation.`n"
$ws.Range("C7").Value = "MSG: None`n`nMSG: The decision-making process has concluded without a choice of movie for Friday. No movie will be selected at this time.`n"
$ws.Range("C8").Value = "MSG: None`n`nMSG: The decision about which movie to show on Friday was not finalized, resulting in no decision being made.`n"
$ws.Range("C9").Value = "MSG: None`n`nMSG: The decision has been recorded as no movie being selected for Friday.`n"
$ws.Range("C10").Value = "MSG: None`n`nMSG: The conversation ended without a decision on which movie to show on Friday.`n"

# Column D (JudgeBotFunctionCalls) update
$ws.Range("D3").Value = "both_movies, "
